# Update gym equipment prices across the three worksheets.
# Prices are stored as literal text (e.g. "$2,072.00"), not numbers, so we
# force text interpretation (NumberFormat "@") while assigning the value,
# then restore the cell's style to "Normal" so no stray number formatting
# sticks around on the cell afterwards.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# Sheet "4x4 Squat Racks"
$ws1 = $wb.Worksheets.Item("4x4 Squat Racks")
Set-TextValue $ws1.Range("C2") '$2,137.00'
Set-TextValue $ws1.Range("C3") '$1,299.99'

# Sheet "Squat Stands"
$ws2 = $wb.Worksheets.Item("Squat Stands")
Set-TextValue $ws2.Range("C2") '$1,545.00'
Set-TextValue $ws2.Range("C3") '$528.99'

# Sheet "Leg Extensions"
$ws3 = $wb.Worksheets.Item("Leg Extensions")
Set-TextValue $ws3.Range("C3") '$2,909.99'
